$wb = $excel.ActiveWorkbook

# --- RUNMANAGER sheet ---
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws1.Range("E2").Value = "1"
$ws1.Range("E3").Value = "2"
$ws1.Range("C4").Value = "no"
$ws1.Range("E2").Select()

# --- DATA sheet ---
$ws2 = $wb.Worksheets.Item("DATA")
$ws2.Range("B3").Value = "no"
$ws2.Range("G2").Value = ""
$ws2.Range("G3").Value = ""
$ws2.Range("G4").Value = ""
$ws2.Range("G5").Value = ""
$ws2.Range("G6").Value = ""
$ws2.Range("C7").Select()
